$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("https://spelevadores.com.br/", "contato@spelevadores.com.br;"),
    @("https://www.elevadoreskorman.com.br/empresas-elevadores-sao-paulo", "vendas@elevadoreskorman.com.br;korman@elevadoreskorman.com.br;comercial@elevadoreskorman.com.br;"),
    @("http://www.emelevadores.com.br/", "contato@emelevadores.com.br;"),
    @("https://coteibem.sindiconet.com.br/fornecedores/manutencao-elevadores/sp/sao-paulo", "contato@coteibem.com.br;"),
    @("http://primac.com.br/", "comercial@primac.com.br;"),
    @("https://retrofitelevadores.com.br/", "contato@elevadoresretrofit.com.br;")
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

$ws.Range("A1:B1").Copy()
$ws.Range("A2:B7").PasteSpecial(-4122)
$ws.Range("A2:B7").Font.Bold = $false
$excel.CutCopyMode = $false
